$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the "Lichtwark"-related columns (B:E) with the values that were
# kept after the deleted columns were removed from the source data.
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 92.73321619510584
$ws.Range("C2").Value = 93.298761155296816
$ws.Range("D2").Value = 93.834676583843432
$ws.Range("E2").Value = 93.780604186449807

$ws.Range("B3").Value = 93.79004036183251
$ws.Range("C3").Value = 96.472162493034546
$ws.Range("D3").Value = 94.534439327533661
$ws.Range("E3").Value = 93.839416201937325

# Update the selection to match the narrower highlighted range left behind
# after the tweak (was B1:AY3, now B1:E3).
$ws.Range("B1:E3").Select()
